# Regenerate save_data column G ("K") values to reflect the new
# computed strike-count (K) metric instead of the old "Strike#" value.
# Only column G (rows 2-18) changes; all other columns are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 0
    4  = 2
    5  = 4
    6  = 2
    7  = 2
    8  = 1
    9  = 2
    10 = 1
    11 = 3
    12 = 4
    13 = 5
    14 = 3
    15 = 4
    16 = 2
    17 = 0
    18 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
